# Sari_Cases_25.xlsx - FluNet template update
# Adds three new "genetic_group" columns, one after each of the three
# "final_res_lineage" / "final_res_lineage2" / "final_res_lineage3" columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert "genetic_group" right after "final_res_lineage" (originally column EZ).
$ws.Columns("EZ").EntireColumn.Insert()
$ws.Range("EZ1").Value = "genetic_group"
$ws.Columns("EZ").ColumnWidth = 13

# 2) Insert "genetic_group2" right after "final_res_lineage2".
#    After the first insertion the target slot is now column FE.
$ws.Columns("FE").EntireColumn.Insert()
$ws.Range("FE1").Value = "genetic_group2"
$ws.Columns("FE").ColumnWidth = 14

# 3) Insert "genetic_group3" right after "final_res_lineage3".
#    After the two prior insertions the target slot is now column FJ.
$ws.Columns("FJ").EntireColumn.Insert()
$ws.Range("FJ1").Value = "genetic_group3"
$ws.Columns("FJ").ColumnWidth = 14
